$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - insert a "condition" column before the existing "date" column,
# shifting the session-2 condition/date pair one column to the right.
$ws.Range("D1").Value = "condition"
$ws.Range("E1").Value = "date"
$ws.Range("F1").Value = "condition"
$ws.Range("G1").Value = "date"

# Row 2 (par id 1 / BA): session 1 condition+date already known; add session 2
# data collection (screen condition + date).
$ws.Range("D2").Value = "VR"
$ws.Range("E2").Value = 20250107
$ws.Range("F2").Value = "screen"
$ws.Range("G2").Value = 20250108

# Row 3 (par id 2 / SY)
$ws.Range("D3").Value = "screen"
$ws.Range("E3").Value = 20250107
$ws.Range("F3").Value = "VR"
$ws.Range("G3").ClearContents()

# Row 4 (par id 3 / NZ)
$ws.Range("D4").Value = "VR"
$ws.Range("E4").Value = 20250108
$ws.Range("F4").Value = "screen"
$ws.Range("G4").ClearContents()

# Row 5 (par id 4 / DT): add "Year since course" value, record session 1 data
# collection (screen condition), and note a TA for the lab course.
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "screen"
$ws.Range("E5").Value = 20250108
$ws.Range("F5").Value = "VR"
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = "TA for the lab course"

# Row 6 (par id 5) - shift condition/date pair from E/G to D/F
$ws.Range("D6").Value = "VR"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "screen"
$ws.Range("G6").ClearContents()

# Row 7 (par id 6) - shift condition/date pair from E/G to D/F
$ws.Range("D7").Value = "screen"
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = "VR"
$ws.Range("G7").ClearContents()

# Selection moves to D5 to match the saved view state.
$ws.Range("D5").Select()
